$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.102.73'
$ws.Range('E2').Value = '  -0.94%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.653.04'
$ws.Range('E3').Value = '  -1.09%  '
$ws.Range('E4').Value = '  -0.55%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.74'
$ws.Range('E5').Value = '  -0.86%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5258'
$ws.Range('E6').Value = '  -1.04%  '
$ws.Range('E7').Value = '  -0.53%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2677'
$ws.Range('E8').Value = '  +0.79%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06373'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.55'
$ws.Range('E10').Value = '  -2.26%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07704'
$ws.Range('E11').Value = '  -1.84%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.599'
$ws.Range('E12').Value = '  +1.24%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.705.29'
$ws.Range('E13').Value = '  +2.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.879.58'
$ws.Range('E14').Value = '  -1.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5628'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8225'
$ws.Range('E16').Value = '  +1.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.53'
$ws.Range('E17').Value = '  -0.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.101.42'
$ws.Range('E18').Value = '  -0.95%  '
$ws.Range('E19').Value = '  -0.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.684'
$ws.Range('E20').Value = '  -0.97%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.36'
$ws.Range('E21').Value = '  +0.51%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '191.23'
$ws.Range('E22').Value = '  -5.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.974'
$ws.Range('E23').Value = '  -1.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.005'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.05'
$ws.Range('E25').Value = '  -0.43%  '
$ws.Range('E26').Value = '  -1.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.249'
$ws.Range('E27').Value = '  -0.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.97'
$ws.Range('E28').Value = '  -1.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.505'
$ws.Range('E29').Value = '  -0.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05645'
$ws.Range('E30').Value = '  -4.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.271'
$ws.Range('E31').Value = '  -1.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.495'
$ws.Range('E32').Value = '  -1.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.382'
$ws.Range('E33').Value = '  +1.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.579'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.795'
$ws.Range('E35').Value = '  -1.37%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9454'
$ws.Range('E36').Value = '  -2.34%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.408'
$ws.Range('E37').Value = '  -1.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5780'
$ws.Range('E38').Value = '  -0.43%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01596'
$ws.Range('E39').Value = '  -1.50%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.972'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8447'
$ws.Range('E41').Value = '  -1.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.004'
$ws.Range('E42').Value = '  -0.59%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.020.73'
$ws.Range('E43').Value = '  -5.43%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.34'
$ws.Range('E44').Value = '  -1.86%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.790.47'
$ws.Range('E45').Value = '  -1.00%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '58.40'
$ws.Range('E46').Value = '  -0.33%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0₈105'
$ws.Range('E47').Value = '  -1.20%  '
$ws.Range('E48').Value = '  +3.76%  '
$ws.Range('B49').Value = 'Frax'
$ws.Range('C49').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.003'
$ws.Range('E49').Value = '  -1.33%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.043'
$ws.Range('E50').Value = '  -0.95%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4342'
$ws.Range('E51').Value = '  -1.73%  '
